# Add 2022-Q4 data:
#  1. Insert a new worksheet "2022-Q4" right after the "总计" (Total) sheet,
#     populated with the new quarter's fund holdings table.
#  2. Insert a new row into "总计" for the 2022-Q4 summary figures, pushing
#     all the older quarters down by one row.
#  3. Keep the workbook's active/selected tab on "2020-Q4" (the last sheet),
#     matching the original file's selection.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Build the new "2022-Q4" sheet.
#    Sheet #4 ("2022-Q1") already has the exact same row/column layout we
#    need (header row + 10 fund rows, columns A:H), so copy it for the
#    formatting/dimensions and then overwrite every value.
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item(4)
$templateSheet.Copy($null, $totalSheet)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

$codes    = @("512980","160629","159855","164818","159805","516620","006048","007413","010677","015675")
$names    = @("广发中证传媒ETF","鹏华中证传媒指数（LOF）A","银华中证影视主题ETF","工银瑞信中证传媒指数（LOF）A","鹏华中证传媒ETF","国泰中证影视主题ETF","长城中证500指数增强A","长城中证500指数增强C","工银瑞信中证传媒指数（LOF）C","鹏华中证传媒指数（LOF）C")
$scales   = @("48.31","7.00","1.01","1.82","1.72","0.71","2.20","0.98","0.27","0.17")
$stockPos = @("99.36","94.29","97.80","93.67","98.37","98.01","94.03","94.03","93.67","94.29")
$posRatio = @("2.85","2.70","6.19","2.69","2.82","6.53","1.18","1.18","2.69","2.70")
$mktVal   = @("1.3768","0.1890","0.0625","0.0490","0.0485","0.0464","0.0260","0.0116","0.0073","0.0046")
$ranks    = @(8,8,4,8,8,4,10,10,8,8)

for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = "'" + $codes[$i]
    $q4.Cells.Item($r, 3).Value = $names[$i]
    $q4.Cells.Item($r, 4).Value = "'" + $scales[$i]
    $q4.Cells.Item($r, 5).Value = "'" + $stockPos[$i]
    $q4.Cells.Item($r, 6).Value = "'" + $posRatio[$i]
    $q4.Cells.Item($r, 7).Value = "'" + $mktVal[$i]
    $q4.Cells.Item($r, 8).Value = $ranks[$i]
}

# ---------------------------------------------------------------------------
# 2. Insert the 2022-Q4 row into "总计", shifting the existing quarters down.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 10
$totalSheet.Cells.Item(2, 4).Value = 1.82

# Column A is a plain 0-based row index; rewrite it for every data row so it
# stays a contiguous sequence after the insert.
$lastRow = $totalSheet.Cells.Item(1,1).End(-4121).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 3. Restore the original active tab (last sheet, "2020-Q4").
# ---------------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Select()
